$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.291.87'
$ws.Range('E2').Value = '  +1.36%  '

$ws.Range('D3').Value = '1.657.69'
$ws.Range('E3').Value = '  +1.18%  '

$ws.Range('E4').Value = '  +1.33%  '

$ws.Range('D5').Value = "'217.77"
$ws.Range('E5').Value = '  +1.08%  '

$ws.Range('E6').Value = '  +1.11%  '

$ws.Range('E7').Value = '  +1.44%  '

$ws.Range('D8').Value = "'0.0641"
$ws.Range('E8').Value = '  +0.24%  '

$ws.Range('E9').Value = '  -0.46%  '

$ws.Range('E10').Value = '  -0.39%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.29"
$ws.Range('E12').Value = '  +0.53%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.661.18'
$ws.Range('E13').Value = '  +1.48%  '

$ws.Range('D14').Value = "'0.545"
$ws.Range('E14').Value = '  +0.06%  '

$ws.Range('D15').Value = "'63.77"
$ws.Range('E15').Value = '  +1.21%  '

$ws.Range('D16').Value = '0.0₃0765'
$ws.Range('E16').Value = '  +0.11%  '

$ws.Range('D17').Value = '26.292.89'
$ws.Range('E17').Value = '  +1.27%  '

$ws.Range('E18').Value = '  +1.35%  '

$ws.Range('D19').Value = "'194.88"
$ws.Range('E19').Value = '  +0.99%  '

$ws.Range('D20').Value = "'4.36"
$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('D21').Value = "'9.82"
$ws.Range('E21').Value = '  -1.05%  '

$ws.Range('D22').Value = "'6.22"
$ws.Range('E22').Value = '  -0.97%  '

$ws.Range('E23').Value = '  +1.70%  '

$ws.Range('D24').Value = "'145.49"
$ws.Range('E24').Value = '  +0.88%  '

$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = "'1.02"
$ws.Range('E25').Value = '  +1.69%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = "'1.79"
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('E27').Value = '  +0.66%  '

$ws.Range('D28').Value = "'15.59"
$ws.Range('E28').Value = '  +0.16%  '

$ws.Range('E29').Value = '  +0.57%  '

$ws.Range('D30').Value = "'0.0491"
$ws.Range('E30').Value = '  -2.48%  '

$ws.Range('E31').Value = '  +1.23%  '

$ws.Range('E32').Value = '  -0.77%  '

$ws.Range('E33').Value = '  +0.30%  '

$ws.Range('E34').Value = '  +1.50%  '

$ws.Range('D35').Value = "'0.907"
$ws.Range('E35').Value = '  +0.57%  '

$ws.Range('D36').Value = '1.142.12'
$ws.Range('E36').Value = '  +0.34%  '

$ws.Range('E37').Value = '  +1.01%  '

$ws.Range('E38').Value = '  -2.14%  '

$ws.Range('E39').Value = '  -0.21%  '

$ws.Range('D40').Value = "'0.804"
$ws.Range('E40').Value = '  +0.70%  '

$ws.Range('D41').Value = "'99.19"
$ws.Range('E41').Value = '  -0.13%  '

$ws.Range('E42').Value = '  -2.42%  '

$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = "'56.69"
$ws.Range('E43').Value = '  +0.04%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = "'1.51"
$ws.Range('E44').Value = '  +1.96%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.0524"
$ws.Range('E45').Value = '  -1.40%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'7.76"
$ws.Range('E46').Value = '  +1.07%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = "'0.419"
$ws.Range('E47').Value = '  +1.07%  '

$ws.Range('B48').Value = 'USDD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D48').Value = "'1.02"
$ws.Range('E48').Value = '  +1.23%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.0945"
$ws.Range('E49').Value = '  -2.04%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = "'1.19"
$ws.Range('E50').Value = '  +2.42%  '

$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = "'5.49"
$ws.Range('E51').Value = '  -0.73%  '
